$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.456.15'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '3.430.28'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''232.72'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").Value = '''620.27'
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("D7").Value = '''1.39'
$ws.Range("E7").Value = '  -2.43%  '
$ws.Range("D8").Value = '''0.395'
$ws.Range("E8").Value = '  -1.68%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '''0.973'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '3.429.86'
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '''42.76'
$ws.Range("E12").Value = '  +3.12%  '
$ws.Range("D13").Value = '''0.199'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").Value = '''6.27'
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '93.337.93'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '4.066.11'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").Value = '''0.0000247'
$ws.Range("E17").Value = '  -1.19%  '
$ws.Range("D18").Value = '''8.20'
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("D19").Value = '3.432.36'
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("E20").Value = '  +4.08%  '
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '''503.05'
$ws.Range("E22").Value = '  +0.98%  '
$ws.Range("D23").Value = '''3.37'
$ws.Range("E23").Value = '  +4.31%  '
$ws.Range("D24").Value = '''0.450'
$ws.Range("E24").Value = '  -4.10%  '
$ws.Range("D25").Value = '''6.65'
$ws.Range("E25").Value = '  +2.57%  '
$ws.Range("E26").Value = '  -3.48%  '
$ws.Range("D27").Value = '''94.99'
$ws.Range("E27").Value = '  +4.38%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = '''11.99'
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.617.78'
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("D30").Value = '''11.46'
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.138'
$ws.Range("E32").Value = '  +1.71%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''2.74'
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("D34").Value = '''0.987'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").Value = '''0.173'
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").Value = '''30.04'
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("D37").Value = '''0.550'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").Value = '''558.85'
$ws.Range("E38").Value = '  +3.20%  '
$ws.Range("D39").Value = '''7.49'
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("D40").Value = '''1.41'
$ws.Range("E40").Value = '  -1.76%  '
$ws.Range("D42").Value = '''0.919'
$ws.Range("E42").Value = '  +1.63%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = '''1.73'
$ws.Range("E44").Value = '  +1.28%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = '''23.69'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("D48").Value = '''0.0410'
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").Value = '''2.12'
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("D51").Value = '''8.08'
$ws.Range("E51").Value = '  +0.61%  '
